# This workbook has 3 sheets:
#   1) "พลศึกษา"   (Physical Education grades)
#   2) "นาฏศิลป์"  (Dance/performing-arts grades)
#   3) "Sheet1"    (lookup table used by defined name grade_edu)
#
# Changes to apply (per the target diff):
#   Sheets 1 & 2:
#     - sheetView: drop the frozen/scrolled topLeftCell, move the
#       selection to Y9
#     - Row 9 height: 135 -> 122.25
#     - Clear the values in S9:U9 (keep their existing style/format)
#   Sheet 3:
#     - sheetView: move the selection from B2:C9 to F20

$wb = $excel.ActiveWorkbook

# --- Sheet 1: พลศึกษา ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Rows.Item(9).RowHeight = 122.25
$ws1.Range("S9:U9").ClearContents()
[void]$ws1.Range("Y9").Select()

# --- Sheet 2: นาฏศิลป์ ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Rows.Item(9).RowHeight = 122.25
$ws2.Range("S9:U9").ClearContents()
[void]$ws2.Range("Y9").Select()

# --- Sheet 3: Sheet1 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
[void]$ws3.Range("F20").Select()

# Restore the originally active/tab-selected sheet
$ws1.Activate()
